$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44313
$ws.Range("M2").Value = 36
$ws.Range("N2").Value = 14000
$ws.Range("O2").Value = 14000
$ws.Range("P2").Value = 14000
$ws.Range("Q2").Value = '$/caja 14 kilos granel'
$ws.Range("R2").Value = 'Provincia de Limarí'
$ws.Range("S2").Value = 1000
$ws.Range("D3").Value = 45050
$ws.Range("L3").Value = 'Especial'
$ws.Range("M3").Value = 56
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 14000
$ws.Range("P3").Value = 14000
$ws.Range("Q3").Value = '$/caja 14 kilos granel'
$ws.Range("T3").Value = 14
$ws.Range("D4").Value = 45050
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("Q4").Value = '$/caja 14 kilos granel'
$ws.Range("S4").Value = 857
$ws.Range("D5").Value = 45040
$ws.Range("L5").Value = 'Especial'
$ws.Range("M5").Value = 65
$ws.Range("N5").Value = 17000
$ws.Range("O5").Value = 17000
$ws.Range("P5").Value = 17000
$ws.Range("Q5").Value = '$/caja 14 kilos granel'
$ws.Range("S5").Value = 1214
$ws.Range("D6").Value = 45040
$ws.Range("M6").Value = 60
$ws.Range("D7").Value = 44259
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("Q7").Value = '$/caja 15 kilos empedrada'
$ws.Range("S7").Value = 800
$ws.Range("T7").Value = 15
$ws.Range("D8").Value = 44252
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 14000
$ws.Range("O8").Value = 14000
$ws.Range("P8").Value = 14000
$ws.Range("Q8").Value = '$/caja 14 kilos empedrada'
$ws.Range("S8").Value = 1000
$ws.Range("D9").Value = 44630
$ws.Range("M9").Value = 75
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("S9").Value = 1071
$ws.Range("D10").Value = 44245
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("Q10").Value = '$/caja 15 kilos granel'
$ws.Range("T10").Value = 15
$ws.Range("D11").Value = 44627
$ws.Range("M11").Value = 56
$ws.Range("N11").Value = 17000
$ws.Range("O11").Value = 17000
$ws.Range("P11").Value = 17000
$ws.Range("Q11").Value = '$/caja 14 kilos empedrada'
$ws.Range("S11").Value = 1214
$ws.Range("D12").Value = 45014
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("Q12").Value = '$/caja 14 kilos empedrada'
$ws.Range("S12").Value = 1071
$ws.Range("D13").Value = 44323
$ws.Range("M13").Value = 60
$ws.Range("D14").Value = 44614
$ws.Range("M14").Value = 54
$ws.Range("D15").Value = 44616
$ws.Range("M15").Value = 70
$ws.Range("N15").Value = 14000
$ws.Range("O15").Value = 14000
$ws.Range("P15").Value = 14000
$ws.Range("Q15").Value = '$/caja 14 kilos empedrada'
$ws.Range("S15").Value = 1000
$ws.Range("T15").Value = 14
$ws.Range("D16").Value = 45006
$ws.Range("M16").Value = 40
$ws.Range("N16").Value = 16000
$ws.Range("O16").Value = 16000
$ws.Range("P16").Value = 16000
$ws.Range("S16").Value = 1143
$ws.Range("D17").Value = 45001
$ws.Range("M17").Value = 50
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("S17").Value = 1143
$ws.Range("D18").Value = 44239
$ws.Range("M18").Value = 70
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 15000
$ws.Range("S18").Value = 1000
$ws.Range("D19").Value = 45042
$ws.Range("L19").Value = 'Especial'
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 17000
$ws.Range("O19").Value = 17000
$ws.Range("P19").Value = 17000
$ws.Range("S19").Value = 1214
$ws.Range("D20").Value = 45042
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 50
$ws.Range("D21").Value = 44278
$ws.Range("M21").Value = 45
$ws.Range("N21").Value = 13000
$ws.Range("O21").Value = 13000
$ws.Range("P21").Value = 13000
$ws.Range("Q21").Value = '$/caja 14 kilos empedrada'
$ws.Range("R21").Value = 'Provincia del Elquí'
$ws.Range("S21").Value = 929
$ws.Range("D22").Value = 45054
$ws.Range("L22").Value = 'Especial'
$ws.Range("M22").Value = 54
$ws.Range("N22").Value = 16000
$ws.Range("O22").Value = 16000
$ws.Range("P22").Value = 16000
$ws.Range("Q22").Value = '$/caja 14 kilos empedrada'
$ws.Range("R22").Value = 'Provincia de Limarí'
$ws.Range("S22").Value = 1143
$ws.Range("D23").Value = 45054
$ws.Range("M23").Value = 50
$ws.Range("D25").Value = 44314
$ws.Range("L25").Value = 'Primera'
$ws.Range("M25").Value = 56
$ws.Range("N25").Value = 14000
$ws.Range("O25").Value = 14000
$ws.Range("P25").Value = 14000
$ws.Range("S25").Value = 1000
$ws.Range("D26").Value = 44312
$ws.Range("M26").Value = 68
$ws.Range("D27").Value = 44270
$ws.Range("M27").Value = 85
$ws.Range("N27").Value = 12000
$ws.Range("O27").Value = 12000
$ws.Range("P27").Value = 12000
$ws.Range("Q27").Value = '$/caja 14 kilos granel'
$ws.Range("R27").Value = 'Provincia del Elquí'
$ws.Range("S27").Value = 857
$ws.Range("T27").Value = 14
$ws.Range("D28").Value = 44592
$ws.Range("M28").Value = 54
$ws.Range("N28").Value = 20000
$ws.Range("O28").Value = 20000
$ws.Range("P28").Value = 20000
$ws.Range("Q28").Value = '$/caja 15 kilos empedrada'
$ws.Range("S28").Value = 1333
$ws.Range("D29").Value = 45044
$ws.Range("M29").Value = 30
$ws.Range("N29").Value = 16000
$ws.Range("O29").Value = 16000
$ws.Range("P29").Value = 16000
$ws.Range("S29").Value = 1143
$ws.Range("D30").Value = 45044
$ws.Range("M30").Value = 30
$ws.Range("D31").Value = 44238
$ws.Range("M31").Value = 60
$ws.Range("N31").Value = 15000
$ws.Range("O31").Value = 15000
$ws.Range("P31").Value = 15000
$ws.Range("Q31").Value = '$/caja 15 kilos granel'
$ws.Range("S31").Value = 1000
$ws.Range("D32").Value = 44315
$ws.Range("M32").Value = 65
$ws.Range("D33").Value = 44322
$ws.Range("M33").Value = 50
$ws.Range("D34").Value = 45043
$ws.Range("L34").Value = 'Especial'
$ws.Range("M34").Value = 45
$ws.Range("N34").Value = 17000
$ws.Range("O34").Value = 17000
$ws.Range("P34").Value = 17000
$ws.Range("Q34").Value = '$/caja 14 kilos granel'
$ws.Range("S34").Value = 1214
$ws.Range("D35").Value = 45043
$ws.Range("M35").Value = 67
$ws.Range("N35").Value = 14000
$ws.Range("O35").Value = 14000
$ws.Range("P35").Value = 14000
$ws.Range("Q35").Value = '$/caja 14 kilos granel'
$ws.Range("S35").Value = 1000
$ws.Range("D36").Value = 44320
$ws.Range("M36").Value = 45
$ws.Range("N36").Value = 14000
$ws.Range("O36").Value = 14000
$ws.Range("P36").Value = 14000
$ws.Range("Q36").Value = '$/caja 14 kilos granel'
$ws.Range("R36").Value = 'Provincia de Limarí'
$ws.Range("S36").Value = 1000
$ws.Range("D37").Value = 44588
$ws.Range("N37").Value = 19000
$ws.Range("O37").Value = 20000
$ws.Range("P37").Value = 19529
$ws.Range("R37").Value = 'Provincia de Limarí'
$ws.Range("S37").Value = 1395
$ws.Range("D38").Value = 44242
$ws.Range("L38").Value = 'Primera'
$ws.Range("N38").Value = 12000
$ws.Range("O38").Value = 12000
$ws.Range("P38").Value = 12000
$ws.Range("Q38").Value = '$/caja 15 kilos granel'
$ws.Range("S38").Value = 800
$ws.Range("T38").Value = 15
$ws.Range("D39").Value = 45015
$ws.Range("M39").Value = 56
$ws.Range("N39").Value = 15000
$ws.Range("O39").Value = 15000
$ws.Range("P39").Value = 15000
$ws.Range("Q39").Value = '$/caja 14 kilos empedrada'
$ws.Range("S39").Value = 1071
$ws.Range("D40").Value = 44271
$ws.Range("M40").Value = 50
$ws.Range("N40").Value = 12000
$ws.Range("O40").Value = 12000
$ws.Range("P40").Value = 12000
$ws.Range("R40").Value = 'Provincia del Elquí'
$ws.Range("S40").Value = 857
$ws.Range("D41").Value = 44260
$ws.Range("M41").Value = 56
$ws.Range("N41").Value = 13000
$ws.Range("O41").Value = 13000
$ws.Range("P41").Value = 13000
